$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1912.7142
$ws.Range("I92").Value = 1399
$ws.Range("J92").Value = 4995
$ws.Range("K92").Value = 1399
$ws.Range("L92").Value = 4995
$ws.Range("M92").Value = -151
$ws.Range("N92").Value = -7491
$ws.Range("H100").Value = 7666.1665
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 7666.1665
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 7666.1665
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -8748.166499999999
$ws.Range("H101").Value = 2643.6875
$ws.Range("I101").Value = 1177.0769
$ws.Range("J101").Value = 8999
$ws.Range("K101").Value = 3531.2307
$ws.Range("L101").Value = 26997
$ws.Range("M101").Value = -1909.2307
$ws.Range("N101").Value = -30241
$ws.Range("H106").Value = 2606785
$ws.Range("I106").Value = 2942827.2
$ws.Range("J106").Value = 2457.75
$ws.Range("K106").Value = 2942827.2
$ws.Range("L106").Value = 2457.75
$ws.Range("M106").Value = -2942196.2
$ws.Range("N106").Value = -3719.75
$ws.Range("H129").Value = 1515.5264
$ws.Range("I129").Value = 821.6875
$ws.Range("K129").Value = 2465.0625
$ws.Range("M129").Value = 2534.9375
$ws.Range("H132").Value = 3195.239
$ws.Range("I132").Value = 2951.9285
$ws.Range("J132").Value = 5750
$ws.Range("K132").Value = 8855.7855
$ws.Range("L132").Value = 17250
$ws.Range("M132").Value = -6325.7855
$ws.Range("N132").Value = -22310
$ws.Range("H137").Value = 2689.2942
$ws.Range("I137").Value = 1798.6666
$ws.Range("J137").Value = 3175.0908
$ws.Range("K137").Value = 5395.9998
$ws.Range("L137").Value = 9525.2724
$ws.Range("M137").Value = -2845.9998
$ws.Range("N137").Value = -14625.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3101.2092
$ws.Range("I32").Value = 2693.1052
$ws.Range("K32").Value = 2693.1052
$ws.Range("M32").Value = -2406.1052
$ws.Range("H74").Value = 1884.8125
$ws.Range("I74").Value = 1849.6923
$ws.Range("J74").Value = 2037
$ws.Range("K74").Value = 1849.6923
$ws.Range("L74").Value = 2037
$ws.Range("M74").Value = -975.6922999999999
$ws.Range("N74").Value = -3785
$ws.Range("H77").Value = 1884.8125
$ws.Range("I77").Value = 1849.6923
$ws.Range("J77").Value = 2037
$ws.Range("K77").Value = 9248.461499999999
$ws.Range("L77").Value = 10185
$ws.Range("M77").Value = -4880.461499999999
$ws.Range("N77").Value = -18921
$ws.Range("H97").Value = 9080.866
$ws.Range("J97").Value = 7168.5
$ws.Range("L97").Value = 7168.5
$ws.Range("N97").Value = -8160.5
$ws.Range("H122").Value = 1979.8695
$ws.Range("I122").Value = 1137.3125
$ws.Range("J122").Value = 3905.7144
$ws.Range("K122").Value = 3411.9375
$ws.Range("L122").Value = 11717.1432
$ws.Range("M122").Value = -961.9375
$ws.Range("N122").Value = -16617.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 44511.668
$ws.Range("I94").Value = 16769.5
$ws.Range("J94").Value = 99996
$ws.Range("K94").Value = 16769.5
$ws.Range("L94").Value = 99996
$ws.Range("M94").Value = -16318.5
$ws.Range("N94").Value = -100898
$ws.Range("H105").Value = 2933.5
$ws.Range("J105").Value = 3498
$ws.Range("L105").Value = 3498
$ws.Range("N105").Value = -6992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2894.6
$ws.Range("I31").Value = 1496.6666
$ws.Range("J31").Value = 3826.5557
$ws.Range("K31").Value = 1496.6666
$ws.Range("L31").Value = 3826.5557
$ws.Range("M31").Value = -1201.6666
$ws.Range("N31").Value = -4416.5557
$ws.Range("H34").Value = 2894.6
$ws.Range("I34").Value = 1496.6666
$ws.Range("J34").Value = 3826.5557
$ws.Range("K34").Value = 1496.6666
$ws.Range("L34").Value = 3826.5557
$ws.Range("M34").Value = -1294.6666
$ws.Range("N34").Value = -4230.5557
$ws.Range("H134").Value = 6245.057
$ws.Range("I134").Value = 5475.433
$ws.Range("J134").Value = 10862.8
$ws.Range("K134").Value = 16426.299
$ws.Range("L134").Value = 32588.4
$ws.Range("M134").Value = -13891.299
$ws.Range("N134").Value = -37658.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 5554.6665
$ws.Range("J118").Value = 8668
$ws.Range("L118").Value = 26004
$ws.Range("N118").Value = -28490
$ws.Range("H120").Value = 15762.25
$ws.Range("I120").Value = 14468
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 43404
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -38566
$ws.Range("N120").Value = -99673
$ws.Range("H133").Value = 23467.916
$ws.Range("I133").Value = 3871.6667
$ws.Range("K133").Value = 11615.0001
$ws.Range("M133").Value = -6555.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8332.388999999999
$ws.Range("I70").Value = 8283.4375
$ws.Range("J70").Value = 8724
$ws.Range("K70").Value = 8283.4375
$ws.Range("L70").Value = 8724
$ws.Range("M70").Value = -8013.4375
$ws.Range("N70").Value = -9264
$ws.Range("H73").Value = 8332.388999999999
$ws.Range("I73").Value = 8283.4375
$ws.Range("J73").Value = 8724
$ws.Range("K73").Value = 8283.4375
$ws.Range("L73").Value = 8724
$ws.Range("M73").Value = -7347.4375
$ws.Range("N73").Value = -10596
$ws.Range("H80").Value = 4444
$ws.Range("I80").Value = 4444
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4444
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3446
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 4444
$ws.Range("I83").Value = 4444
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 22220
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -17228
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 81349.2
$ws.Range("J36").Value = 81349.2
$ws.Range("L36").Value = 81349.2
$ws.Range("N36").Value = -82473.2
$ws.Range("H68").Value = 2381.8948
$ws.Range("I68").Value = 2163.3635
$ws.Range("K68").Value = 2163.3635
$ws.Range("M68").Value = -1414.3635
$ws.Range("H71").Value = 2381.8948
$ws.Range("I71").Value = 2163.3635
$ws.Range("K71").Value = 10816.8175
$ws.Range("M71").Value = -7072.817499999999
$ws.Range("H93").Value = 1254.8
$ws.Range("I93").Value = 1048.25
$ws.Range("K93").Value = 1048.25
$ws.Range("M93").Value = 199.75
$ws.Range("H122").Value = 7467.15
$ws.Range("I122").Value = 6949.5
$ws.Range("J122").Value = 7596.5625
$ws.Range("K122").Value = 20848.5
$ws.Range("L122").Value = 22789.6875
$ws.Range("M122").Value = -18398.5
$ws.Range("N122").Value = -27689.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1597.4445
$ws.Range("I100").Value = 624.8570999999999
$ws.Range("K100").Value = 1249.7142
$ws.Range("M100").Value = -708.7141999999999
